# Update cryptos list values per data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.603.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5337"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("E8").Value = "  +5.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07807"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.120"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.337"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.594"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.58%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.856.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.67%  "
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06562"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.55%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.096"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.614.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.234"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.043.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.416"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.150"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1120"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.757"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.652"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07345"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2270"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02356"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.53%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.952"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.234"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6312"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.201"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.90%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.395"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5951"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.709"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.194"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06957"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.07%  "
